$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data columns so numeric-looking strings
# (e.g. "1.004", "41.39") are preserved as text, matching the source data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.874.92"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "1.782.79"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "311.04"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5112"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "0.3761"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").Value = "0.07765"
$ws.Range("E9").Value = "  -8.40%  "
$ws.Range("D10").Value = "41.39"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "1.083"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "6.179"
$ws.Range("E13").Value = "  -3.97%  "
$ws.Range("D14").Value = "20.14"
$ws.Range("E14").Value = "  -4.28%  "
$ws.Range("D15").Value = "1.782.56"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "7.175"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "91.83"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "0.00001068"
$ws.Range("E18").Value = "  -6.22%  "
$ws.Range("D19").Value = "0.06528"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "17.00"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").Value = "5.901"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "27.940.77"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "158.05"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").Value = "1.986.97"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "2.350"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").Value = "122.16"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("D31").Value = "0.1072"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "1.036"
$ws.Range("E32").Value = "  -5.07%  "
$ws.Range("D33").Value = "3.629"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "5.472"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "0.07085"
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("D36").Value = "0.02303"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").Value = "0.2118"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").Value = "8.545"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.998"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "11.46"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").Value = "0.6084"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.152"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.332"
$ws.Range("E44").Value = "  -5.11%  "
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5937"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.723"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "126.27"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.210"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.890"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06721"
$ws.Range("E51").Value = "  -3.77%  "
